$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$dateParagraph = $d.Paragraphs.Item(1)
$dateParagraph.Range.Text = "2024-02-29 Thursday"

# Update each arithmetic-expression cell in the practice table, addressed by
# (row, column) position so duplicate expression text across cells does not
# cause cross-cell replacement collisions.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "70+21="
$t.Cell(1, 2).Range.Text = "16+16="
$t.Cell(1, 3).Range.Text = "60+32="
$t.Cell(1, 4).Range.Text = "69-36="
$t.Cell(1, 5).Range.Text = "11+58="
$t.Cell(2, 1).Range.Text = "31-7="
$t.Cell(2, 2).Range.Text = "4+45="
$t.Cell(2, 3).Range.Text = "20+53="
$t.Cell(2, 4).Range.Text = "67-58="
$t.Cell(2, 5).Range.Text = "44+51="
$t.Cell(3, 1).Range.Text = "86-79="
$t.Cell(3, 2).Range.Text = "52-23="
$t.Cell(3, 3).Range.Text = "70-54="
$t.Cell(3, 4).Range.Text = "75-21="
$t.Cell(3, 5).Range.Text = "78-11="
$t.Cell(4, 1).Range.Text = "46-8="
$t.Cell(4, 2).Range.Text = "9+62="
$t.Cell(4, 3).Range.Text = "65-37="
$t.Cell(4, 4).Range.Text = "10+79="
$t.Cell(4, 5).Range.Text = "92-6="
$t.Cell(5, 1).Range.Text = "47+23="
$t.Cell(5, 2).Range.Text = "16+35="
$t.Cell(5, 3).Range.Text = "68-58="
$t.Cell(5, 4).Range.Text = "94-41="
$t.Cell(5, 5).Range.Text = "78-77="
$t.Cell(6, 1).Range.Text = "35+25="
$t.Cell(6, 2).Range.Text = "46-38="
$t.Cell(6, 3).Range.Text = "20+1="
$t.Cell(6, 4).Range.Text = "15-4="
$t.Cell(6, 5).Range.Text = "54-10="
$t.Cell(7, 1).Range.Text = "89+1="
$t.Cell(7, 2).Range.Text = "79+6="
$t.Cell(7, 3).Range.Text = "28+63="
$t.Cell(7, 4).Range.Text = "90-15="
$t.Cell(7, 5).Range.Text = "14-12="
$t.Cell(8, 1).Range.Text = "58+6="
$t.Cell(8, 2).Range.Text = "37+32="
$t.Cell(8, 3).Range.Text = "6+34="
$t.Cell(8, 4).Range.Text = "64-43="
$t.Cell(8, 5).Range.Text = "42-23="
$t.Cell(9, 1).Range.Text = "55-30="
$t.Cell(9, 2).Range.Text = "88-35="
$t.Cell(9, 3).Range.Text = "5+3="
$t.Cell(9, 4).Range.Text = "41+53="
$t.Cell(9, 5).Range.Text = "45+10="
$t.Cell(10, 1).Range.Text = "12+46="
$t.Cell(10, 2).Range.Text = "77-47="
$t.Cell(10, 3).Range.Text = "96-61="
$t.Cell(10, 4).Range.Text = "46-44="
$t.Cell(10, 5).Range.Text = "76-1="
$t.Cell(11, 1).Range.Text = "42-3="
$t.Cell(11, 2).Range.Text = "56+39="
$t.Cell(11, 3).Range.Text = "54-26="
$t.Cell(11, 4).Range.Text = "26+42="
$t.Cell(11, 5).Range.Text = "3+60="
$t.Cell(12, 1).Range.Text = "28-16="
$t.Cell(12, 2).Range.Text = "90-37="
$t.Cell(12, 3).Range.Text = "64-11="
$t.Cell(12, 4).Range.Text = "52+40="
$t.Cell(12, 5).Range.Text = "95-68="
$t.Cell(13, 1).Range.Text = "69+14="
$t.Cell(13, 2).Range.Text = "58-3="
$t.Cell(13, 3).Range.Text = "79+9="
$t.Cell(13, 4).Range.Text = "63-38="
$t.Cell(13, 5).Range.Text = "91-2="
$t.Cell(14, 1).Range.Text = "65-23="
$t.Cell(14, 2).Range.Text = "74-42="
$t.Cell(14, 3).Range.Text = "70+14="
$t.Cell(14, 4).Range.Text = "72-65="
$t.Cell(14, 5).Range.Text = "38+0="
$t.Cell(15, 1).Range.Text = "78+21="
$t.Cell(15, 2).Range.Text = "28-23="
$t.Cell(15, 3).Range.Text = "18+35="
$t.Cell(15, 4).Range.Text = "85-63="
$t.Cell(15, 5).Range.Text = "19+76="
$t.Cell(16, 1).Range.Text = "53-24="
$t.Cell(16, 2).Range.Text = "77-1="
$t.Cell(16, 3).Range.Text = "2+32="
$t.Cell(16, 4).Range.Text = "62-53="
$t.Cell(16, 5).Range.Text = "48-3="
$t.Cell(17, 1).Range.Text = "38-29="
$t.Cell(17, 2).Range.Text = "82-57="
$t.Cell(17, 3).Range.Text = "68-29="
$t.Cell(17, 4).Range.Text = "31+37="
$t.Cell(17, 5).Range.Text = "36-4="
$t.Cell(18, 1).Range.Text = "80-9="
$t.Cell(18, 2).Range.Text = "33-30="
$t.Cell(18, 3).Range.Text = "27+32="
$t.Cell(18, 4).Range.Text = "29+44="
$t.Cell(18, 5).Range.Text = "52+46="
$t.Cell(19, 1).Range.Text = "39-21="
$t.Cell(19, 2).Range.Text = "15+20="
$t.Cell(19, 3).Range.Text = "16+82="
$t.Cell(19, 4).Range.Text = "60-6="
$t.Cell(19, 5).Range.Text = "57+25="
$t.Cell(20, 1).Range.Text = "24-1="
$t.Cell(20, 2).Range.Text = "95-50="
$t.Cell(20, 3).Range.Text = "53-9="
$t.Cell(20, 4).Range.Text = "65+11="
$t.Cell(20, 5).Range.Text = "59+2="
